$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "K" = 0
    "L" = 1
    "O" = 13.250951766967773
    "P" = 38.024566650390625
    "Q" = 29.036596298217773
    "R" = 8.9879703521728516
    "S" = 3.2160265445709229
    "T" = 35.548210144042969
    "U" = 8.6075477600097656
    "V" = 26.940662384033203
    "W" = 1
    "Y" = 17.446767807006836
    "Z" = 21.880966186523438
    "AA" = 1.1773288249969482
    "AB" = 20.703638076782227
    "AC" = 8.8961515426635742
    "AD" = 2.1978754997253418
    "AE" = 2.1978754997253418
    "AF" = 0
    "AH" = 0
    "AI" = 7.0325741767883301
    "AJ" = 29.231773376464844
    "AK" = 29.231773376464844
    "AM" = 4.6432280540466309
    "AN" = 8.4582910537719727
    "AO" = 8.4582910537719727
    "AQ" = 1
    "AS" = 2.3537311553955078
    "AT" = 32.880050659179688
    "AU" = 9.2514591217041016
    "AV" = 23.628591537475586
    "AW" = 8.6124334335327148
    "AX" = 13.234278678894043
    "AY" = 13.234278678894043
    "AZ" = 0
    "BC" = 1.9467545747756958
    "BD" = 44.056102752685547
    "BE" = 39.308578491210938
    "BF" = 4.7475242614746094
    "BG" = 3.0712547302246094
    "BH" = 47.781539916992188
    "BI" = 43.099933624267578
    "BJ" = 4.6816062927246094
    "BK" = 1
    "BM" = 18.21574592590332
    "BN" = 42.9378662109375
    "BO" = 7.8077750205993652
    "BP" = 35.130092620849609
    "BQ" = 14.384234428405762
    "BR" = 46.834102630615234
    "BS" = 3.9423618316650391
    "BT" = 42.891738891601563
    "BU" = 10.666536331176758
    "BV" = 7.4639077186584473

}

foreach ($row in 2,3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

